# Auto-generated Excel COM-interop script applying the Aegis_Profits market-data refresh.
# For each (sheet, row) the currentAveragePrice* (H/I/J), LevePrice* (K/L), and LeveProfit* (M/N)
# columns are updated to reflect the latest market snapshot pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 727490.25
$ws.Range("I58").Value = 1508709.1
$ws.Range("J58").Value = 2072.6428
$ws.Range("K58").Value = 4526127.300000001
$ws.Range("L58").Value = 6217.928400000001
$ws.Range("M58").Value = -4525977.300000001
$ws.Range("N58").Value = -6517.928400000001

$ws.Range("H69").Value = 3507.182
$ws.Range("J69").Value = 3175.4443
$ws.Range("L69").Value = 9526.332900000001
$ws.Range("N69").Value = -11274.3329

$ws.Range("H72").Value = 3507.182
$ws.Range("J72").Value = 3175.4443
$ws.Range("L72").Value = 28578.9987
$ws.Range("N72").Value = -37314.9987

$ws.Range("H100").Value = 1437.5
$ws.Range("I100").Value = 983.3333
$ws.Range("J100").Value = 2800
$ws.Range("K100").Value = 983.3333
$ws.Range("L100").Value = 2800
$ws.Range("M100").Value = -442.3333
$ws.Range("N100").Value = -3882

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 24900
$ws.Range("J24").Value = 24900
$ws.Range("L24").Value = 24900
$ws.Range("N24").Value = -25648

$ws.Range("H61").Value = 1499.0167
$ws.Range("I61").Value = 1030.7872
$ws.Range("J61").Value = 3191.8462
$ws.Range("K61").Value = 1030.7872
$ws.Range("L61").Value = 3191.8462
$ws.Range("M61").Value = -818.7872
$ws.Range("N61").Value = -3615.8462

$ws.Range("H74").Value = 736.2
$ws.Range("I74").Value = 695.875
$ws.Range("K74").Value = 695.875
$ws.Range("M74").Value = 178.125

$ws.Range("H77").Value = 736.2
$ws.Range("I77").Value = 695.875
$ws.Range("K77").Value = 3479.375
$ws.Range("M77").Value = 888.625

$ws.Range("H92").Value = 7350
$ws.Range("J92").Value = 7350
$ws.Range("L92").Value = 7350
$ws.Range("N92").Value = -12342

$ws.Range("H100").Value = 24900
$ws.Range("J100").Value = 24900
$ws.Range("L100").Value = 24900
$ws.Range("N100").Value = -27064

$ws.Range("H110").Value = 38543016
$ws.Range("I110").Value = 43570268
$ws.Range("J110").Value = 741
$ws.Range("K110").Value = 43570268
$ws.Range("L110").Value = 741
$ws.Range("M110").Value = -43568223
$ws.Range("N110").Value = -4831

$ws.Range("H112").Value = 11071.75
$ws.Range("J112").Value = 11071.75
$ws.Range("L112").Value = 11071.75
$ws.Range("N112").Value = -14025.75

$ws.Range("H114").Value = 28000
$ws.Range("J114").Value = 28000
$ws.Range("L114").Value = 28000
$ws.Range("N114").Value = -36678

$ws.Range("H136").Value = 1499.0167
$ws.Range("I136").Value = 1030.7872
$ws.Range("J136").Value = 3191.8462
$ws.Range("K136").Value = 3092.3616
$ws.Range("L136").Value = 9575.5386
$ws.Range("M136").Value = -542.3616000000002
$ws.Range("N136").Value = -14675.5386

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 92652.27
$ws.Range("I105").Value = 64282.438
$ws.Range("J105").Value = 168305.17
$ws.Range("K105").Value = 64282.438
$ws.Range("L105").Value = 168305.17
$ws.Range("M105").Value = -62535.438
$ws.Range("N105").Value = -171799.17

$ws.Range("H110").Value = 34850
$ws.Range("J110").Value = 34850
$ws.Range("L110").Value = 34850
$ws.Range("N110").Value = -43030

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1704.7407
$ws.Range("I58").Value = 1459.7273
$ws.Range("J58").Value = 2782.8
$ws.Range("K58").Value = 1459.7273
$ws.Range("L58").Value = 2782.8
$ws.Range("M58").Value = -1256.7273
$ws.Range("N58").Value = -3188.8

$ws.Range("H106").Value = 26300
$ws.Range("J106").Value = 26300
$ws.Range("L106").Value = 26300
$ws.Range("N106").Value = -28824

$ws.Range("H136").Value = 1704.7407
$ws.Range("I136").Value = 1459.7273
$ws.Range("J136").Value = 2782.8
$ws.Range("K136").Value = 4379.1819
$ws.Range("L136").Value = 8348.400000000001
$ws.Range("M136").Value = -1829.1819
$ws.Range("N136").Value = -13448.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1139.6976
$ws.Range("I5").Value = 1054.8667
$ws.Range("J5").Value = 1185.1428
$ws.Range("K5").Value = 3164.6001
$ws.Range("L5").Value = 3555.4284
$ws.Range("M5").Value = -3052.6001
$ws.Range("N5").Value = -3779.4284

$ws.Range("H135").Value = 1139.6976
$ws.Range("I135").Value = 1054.8667
$ws.Range("J135").Value = 1185.1428
$ws.Range("K135").Value = 9493.800300000001
$ws.Range("L135").Value = 10666.2852
$ws.Range("M135").Value = -6958.800300000001
$ws.Range("N135").Value = -15736.2852

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 77004140
$ws.Range("I80").Value = 100104780
$ws.Range("J80").Value = 2000
$ws.Range("K80").Value = 100104780
$ws.Range("L80").Value = 2000
$ws.Range("M80").Value = -100103782
$ws.Range("N80").Value = -3996

$ws.Range("H83").Value = 77004140
$ws.Range("I83").Value = 100104780
$ws.Range("J83").Value = 2000
$ws.Range("K83").Value = 500523900
$ws.Range("L83").Value = 10000
$ws.Range("M83").Value = -500518908
$ws.Range("N83").Value = -19984

$ws.Range("H97").Value = 52633936
$ws.Range("I97").Value = 62502190
$ws.Range("J97").Value = 3240.6667
$ws.Range("K97").Value = 62502190
$ws.Range("L97").Value = 3240.6667
$ws.Range("M97").Value = -62501694
$ws.Range("N97").Value = -4232.6667

$ws.Range("H113").Value = 1880.6666
$ws.Range("I113").Value = 1900
$ws.Range("J113").Value = 1876.8
$ws.Range("K113").Value = 1900
$ws.Range("L113").Value = 1876.8
$ws.Range("M113").Value = 270
$ws.Range("N113").Value = -6216.8

$ws.Range("H122").Value = 1651.1875
$ws.Range("I122").Value = 1448.9474
$ws.Range("K122").Value = 4346.8422
$ws.Range("M122").Value = -1896.8422

$ws.Range("H123").Value = 28771.428
$ws.Range("J123").Value = 28771.428
$ws.Range("L123").Value = 28771.428
$ws.Range("N123").Value = -33671.428

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 511525.1
$ws.Range("I46").Value = 331.66666
$ws.Range("J46").Value = 730608
$ws.Range("K46").Value = 331.66666
$ws.Range("L46").Value = 730608
$ws.Range("M46").Value = -143.66666
$ws.Range("N46").Value = -730984

$ws.Range("H82").Value = 2015.2222
$ws.Range("I82").Value = 1732.8572
$ws.Range("J82").Value = 2194.9092
$ws.Range("K82").Value = 1732.8572
$ws.Range("L82").Value = 2194.9092
$ws.Range("M82").Value = -1371.8572
$ws.Range("N82").Value = -2916.9092

$ws.Range("H85").Value = 2015.2222
$ws.Range("I85").Value = 1732.8572
$ws.Range("J85").Value = 2194.9092
$ws.Range("K85").Value = 1732.8572
$ws.Range("L85").Value = 2194.9092
$ws.Range("M85").Value = -484.8571999999999
$ws.Range("N85").Value = -4690.9092

$ws.Range("H94").Value = 28188.334
$ws.Range("J94").Value = 28188.334
$ws.Range("L94").Value = 28188.334
$ws.Range("N94").Value = -29540.334

$ws.Range("H110").Value = 29333.334
$ws.Range("J110").Value = 29333.334
$ws.Range("L110").Value = 29333.334
$ws.Range("N110").Value = -37513.334

$ws.Range("H136").Value = 1771.6285
$ws.Range("I136").Value = 1549.625
$ws.Range("J136").Value = 2256
$ws.Range("K136").Value = 4648.875
$ws.Range("L136").Value = 6768
$ws.Range("M136").Value = -2098.875
$ws.Range("N136").Value = -11868

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 10360
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 10360
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 10360
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -11342

$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H103").Value = 24250
$ws.Range("J103").Value = 24250
$ws.Range("L103").Value = 24250
$ws.Range("N103").Value = -26594

$ws.Range("H104").Value = 24000
$ws.Range("J104").Value = 24000
$ws.Range("L104").Value = 24000
$ws.Range("N104").Value = -30988

$ws.Range("H105").Value = 42125
$ws.Range("J105").Value = 42125
$ws.Range("L105").Value = 42125
$ws.Range("N105").Value = -49113
